$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 556.36365
$ws.Range("I107").Value = 512
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 512
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1408
$ws.Range("N107").Value = -4840
$ws.Range("H137").Value = 1415.4073
$ws.Range("I137").Value = 1261.7778
$ws.Range("J137").Value = 1722.6666
$ws.Range("K137").Value = 3785.3334
$ws.Range("L137").Value = 5167.9998
$ws.Range("M137").Value = -1235.3334
$ws.Range("N137").Value = -10267.9998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2447.7273
$ws.Range("I61").Value = 2045.8334
$ws.Range("J61").Value = 2930
$ws.Range("K61").Value = 2045.8334
$ws.Range("L61").Value = 2930
$ws.Range("M61").Value = -1833.8334
$ws.Range("N61").Value = -3354
$ws.Range("H63").Value = 836157.9399999999
$ws.Range("I63").Value = 911963.2
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 911963.2
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -911277.2
$ws.Range("N63").Value = -3672
$ws.Range("H66").Value = 836157.9399999999
$ws.Range("I66").Value = 911963.2
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 4559816
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -4556384
$ws.Range("N66").Value = -18364
$ws.Range("H74").Value = 1938.3846
$ws.Range("I74").Value = 2188.7778
$ws.Range("J74").Value = 1375
$ws.Range("K74").Value = 2188.7778
$ws.Range("L74").Value = 1375
$ws.Range("M74").Value = -1314.7778
$ws.Range("N74").Value = -3123
$ws.Range("H76").Value = 34000
$ws.Range("I76").Value = 34000
$ws.Range("K76").Value = 34000
$ws.Range("M76").Value = -33662
$ws.Range("H77").Value = 1938.3846
$ws.Range("I77").Value = 2188.7778
$ws.Range("J77").Value = 1375
$ws.Range("K77").Value = 10943.889
$ws.Range("L77").Value = 6875
$ws.Range("M77").Value = -6575.888999999999
$ws.Range("N77").Value = -15611
$ws.Range("H79").Value = 34000
$ws.Range("I79").Value = 34000
$ws.Range("K79").Value = 34000
$ws.Range("M79").Value = -32830
$ws.Range("H122").Value = 1893.4445
$ws.Range("I122").Value = 1698.1428
$ws.Range("J122").Value = 2103.7693
$ws.Range("K122").Value = 5094.428400000001
$ws.Range("L122").Value = 6311.3079
$ws.Range("M122").Value = -2644.428400000001
$ws.Range("N122").Value = -11211.3079
$ws.Range("H136").Value = 2447.7273
$ws.Range("I136").Value = 2045.8334
$ws.Range("J136").Value = 2930
$ws.Range("K136").Value = 6137.5002
$ws.Range("L136").Value = 8790
$ws.Range("M136").Value = -3587.5002
$ws.Range("N136").Value = -13890
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 73501.71000000001
$ws.Range("I134").Value = 126378
$ws.Range("K134").Value = 379134
$ws.Range("M134").Value = -376599
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1040.5294
$ws.Range("I16").Value = 830
$ws.Range("J16").Value = 1724.75
$ws.Range("K16").Value = 830
$ws.Range("L16").Value = 1724.75
$ws.Range("M16").Value = -543
$ws.Range("N16").Value = -2298.75
$ws.Range("H31").Value = 4169073.2
$ws.Range("I31").Value = 2991.9092
$ws.Range("J31").Value = 7694219
$ws.Range("K31").Value = 2991.9092
$ws.Range("L31").Value = 7694219
$ws.Range("M31").Value = -2696.9092
$ws.Range("N31").Value = -7694809
$ws.Range("H34").Value = 4169073.2
$ws.Range("I34").Value = 2991.9092
$ws.Range("J34").Value = 7694219
$ws.Range("K34").Value = 2991.9092
$ws.Range("L34").Value = 7694219
$ws.Range("M34").Value = -2789.9092
$ws.Range("N34").Value = -7694623
$ws.Range("H58").Value = 1247
$ws.Range("I58").Value = 1438.875
$ws.Range("J58").Value = 940
$ws.Range("K58").Value = 1438.875
$ws.Range("L58").Value = 940
$ws.Range("M58").Value = -1235.875
$ws.Range("N58").Value = -1346
$ws.Range("H74").Value = 34340
$ws.Range("J74").Value = 34340
$ws.Range("L74").Value = 34340
$ws.Range("N74").Value = -36088
$ws.Range("H77").Value = 34340
$ws.Range("J77").Value = 34340
$ws.Range("L77").Value = 103020
$ws.Range("N77").Value = -111756
$ws.Range("H113").Value = 1040.5294
$ws.Range("I113").Value = 830
$ws.Range("J113").Value = 1724.75
$ws.Range("K113").Value = 830
$ws.Range("L113").Value = 1724.75
$ws.Range("M113").Value = 1340
$ws.Range("N113").Value = -6064.75
$ws.Range("H122").Value = 739.0476
$ws.Range("I122").Value = 754.7692
$ws.Range("J122").Value = 713.5
$ws.Range("K122").Value = 2264.3076
$ws.Range("L122").Value = 2140.5
$ws.Range("M122").Value = 185.6923999999999
$ws.Range("N122").Value = -7040.5
$ws.Range("H132").Value = 2575.4546
$ws.Range("I132").Value = 2181.5908
$ws.Range("J132").Value = 3363.182
$ws.Range("K132").Value = 6544.7724
$ws.Range("L132").Value = 10089.546
$ws.Range("M132").Value = -4014.7724
$ws.Range("N132").Value = -15149.546
$ws.Range("H134").Value = 1126.8948
$ws.Range("I134").Value = 1099.8235
$ws.Range("K134").Value = 3299.4705
$ws.Range("M134").Value = -764.4704999999999
$ws.Range("H136").Value = 1247
$ws.Range("I136").Value = 1438.875
$ws.Range("J136").Value = 940
$ws.Range("K136").Value = 4316.625
$ws.Range("L136").Value = 2820
$ws.Range("M136").Value = -1766.625
$ws.Range("N136").Value = -7920
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 670.7895
$ws.Range("I107").Value = 597.6667
$ws.Range("K107").Value = 597.6667
$ws.Range("M107").Value = 1322.3333
$ws.Range("H122").Value = 2201.1177
$ws.Range("I122").Value = 1302.2
$ws.Range("J122").Value = 2575.6667
$ws.Range("K122").Value = 3906.6
$ws.Range("L122").Value = 7727.000100000001
$ws.Range("M122").Value = -1456.6
$ws.Range("N122").Value = -12627.0001
$ws.Range("H126").Value = 1670571.5
$ws.Range("I126").Value = 5258.5713
$ws.Range("J126").Value = 2567278.5
$ws.Range("K126").Value = 15775.7139
$ws.Range("L126").Value = 7701835.5
$ws.Range("M126").Value = -13305.7139
$ws.Range("N126").Value = -7706775.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3373.75
$ws.Range("I61").Value = 2250
$ws.Range("J61").Value = 4497.5
$ws.Range("K61").Value = 2250
$ws.Range("L61").Value = 4497.5
$ws.Range("M61").Value = -2048
$ws.Range("N61").Value = -4901.5
$ws.Range("H113").Value = 3373.75
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 4497.5
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 4497.5
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -8837.5
$ws.Range("H122").Value = 2071.5715
$ws.Range("I122").Value = 2790.4
$ws.Range("J122").Value = 1672.2222
$ws.Range("K122").Value = 8371.200000000001
$ws.Range("L122").Value = 5016.6666
$ws.Range("M122").Value = -5921.200000000001
$ws.Range("N122").Value = -9916.6666
$ws.Range("H132").Value = 14527.852
$ws.Range("J132").Value = 4057.8462
$ws.Range("L132").Value = 12173.5386
$ws.Range("N132").Value = -17233.5386
$ws.Range("H136").Value = 4831.567
$ws.Range("I136").Value = 7683.1333
$ws.Range("K136").Value = 23049.3999
$ws.Range("M136").Value = -20499.3999
$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 35000
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 667.3333
$ws.Range("I107").Value = 667.3333
$ws.Range("K107").Value = 2001.9999
$ws.Range("M107").Value = -81.99990000000003
$ws.Range("H122").Value = 1163.15
$ws.Range("I122").Value = 1025.7273
$ws.Range("J122").Value = 1331.1111
$ws.Range("K122").Value = 3077.1819
$ws.Range("L122").Value = 3993.3333
$ws.Range("M122").Value = -627.1819
$ws.Range("N122").Value = -8893.3333
$ws.Range("H132").Value = 2372.3333
$ws.Range("I132").Value = 1598.375
$ws.Range("K132").Value = 4795.125
$ws.Range("M132").Value = -2265.125
$ws.Range("H136").Value = 5635.96
$ws.Range("J136").Value = 865.25
$ws.Range("L136").Value = 2595.75
$ws.Range("N136").Value = -7695.75
